$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit "Restored from revision ... Type: SAVE" re-saves the workbook; the
# only actual data change in the diff is cell C10 (row 10, the "R30" rule's
# "From" value) going from 18 to 1.
$ws.Range("C10").Value = 1
